$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), shifting
# existing data rows 2-22 down to 4-24.
$ws.Rows.Item(2).Resize(2).Insert() | Out-Null

# Populate the two newly inserted rows with the new data points.
$newRows = @(
    @(-0.1360424668951469, 0.3186911859295585, 0.05561650341207323),
    @(-0.2708076590841467, 0.3116106986999512, 0.0008329986171288634)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# The last three original rows (old 20, 21, 22 -> now shifted to 22, 23, 24)
# are no longer present in the final data; delete them.
$ws.Rows.Item(22).Resize(3).Delete() | Out-Null

$wb.Save()
